$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '63.899.81'
$c.Style = "Normal"
$ws.Range('E2').Value = '  -0.27%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.751.98'
$c.Style = "Normal"
$ws.Range('E3').Value = '  +0.63%  '
$ws.Range('E4').Value = '  +0.02%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '580.01'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -1.98%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '159.82'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +4.16%  '
$ws.Range('E7').Value = '  +0.17%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.613'
$c.Style = "Normal"
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('E9').Value = '  -0.70%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '5.90'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -12.19%  '
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('E12').Value = '  +0.07%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '3.240.04'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('E14').Value = '  +1.52%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '63.867.31'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -0.06%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '0.0000155'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +1.80%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '2.757.37'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -0.16%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '12.32'
$c.Style = "Normal"
$ws.Range('E18').Value = '  +1.90%  '
$ws.Range('E19').Value = '  +1.41%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '362.67'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -0.70%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '6.89'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -1.83%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.572'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +6.14%  '
$ws.Range('E23').Value = '  +0.54%  '
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('E25').Value = '  +2.79%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '8.68'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +0.18%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +0.30%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '0.0₃0940'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +3.18%  '
$ws.Range('E29').Value = '  -1.39%  '
$ws.Range('E30').Value = '  +0.27%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '1.26'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +4.69%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '168.70'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -2.40%  '
$ws.Range('E33').Value = '  +0.06%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '20.59'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -0.14%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '5.02'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +3.56%  '
$ws.Range('E36').Value = '  +2.76%  '
$ws.Range('E37').Value = '  +1.86%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('E39').Value = '  +0.02%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '6.16'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +10.03%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '333.78'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -4.25%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '39.54'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +1.46%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '22.14'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('E44').Value = '  +1.36%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.0259'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '0.641'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -1.20%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '136.81'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -4.79%  '
$ws.Range('E49').Value = '  +0.51%  '
$ws.Range('E50').Value = '  +0.22%  '
$ws.Range('E51').Value = '  +0.70%  '
